$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix "Voltage for lod cell is 2.5V, so not bad!!"
#    -> "Measured voltage for load cell is 2.506V, so not bad!!"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Voltage for lod cell is 2.5V, so not bad!!", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Measured voltage for load cell is 2.506V, so not bad!!", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Fix the typo "toinight" -> "tonight"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Weird behavior toinight:", $false, $false, $false, $false, $false,
    $true, 1, $false, "Weird behavior tonight:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Fix "If I print value and Last reding, the delay is 16us"
#    -> "If I print value and last reading, the delay is 16us"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "If I print value and Last reding, the delay is 16us", $false, $false,
    $false, $false, $false, $true, 1, $false,
    "If I print value and last reading, the delay is 16us", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Merge "Uhhh" + (bookmark) + "lala" into a single run "Uhhhlala" and
#    drop the old _GoBack bookmark from that spot (it gets re-added at the
#    end of the document below).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Content.Find.Execute(
    "Uhhhlala", $false, $false, $false, $false, $false, $true, 1, $false,
    "Uhhhlala", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Append a new blank paragraph followed by a paragraph with the new
#    "ALL NEW: ..." note after the "Uhhhlala" paragraph.
# ---------------------------------------------------------------------------
$uhhhlalaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($uhhhlalaIndex)
$endOfDoc = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endOfDoc.InsertBefore("`r")

$blankPara = $d.Paragraphs.Item($uhhhlalaIndex + 1)
$afterBlank = $blankPara.Range
$afterBlank.Collapse(0)
$afterBlank.InsertAfter("`r")

$notePara = $d.Paragraphs.Item($uhhhlalaIndex + 2)
$noteBody = $d.Range($notePara.Range.Start, $notePara.Range.End - 1)
$noteBody.Text = "ALL NEW: abandoned the RMT approach and now everything is read-out in interrupt, works nicely. The scale is accurate to 0.1g which is quite nice!!"

# ---------------------------------------------------------------------------
# 6. Re-create the _GoBack bookmark at the end of the new note paragraph
#    (collapsed, right after the text and before the paragraph mark).
#    A bookmark cannot be Add()-ed directly at a paragraph-mark position, so
#    a one-character placeholder is used to hold the spot and then removed.
# ---------------------------------------------------------------------------
$notePara = $d.Paragraphs.Item($uhhhlalaIndex + 2)
$endPos = $notePara.Range.End - 1
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertBefore("X")
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$d.Range($endPos, $endPos + 1).Delete()
